$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Supports tax on world top 1% to finance global poverty reduction`n(Additional 15% tax on income over [`$120k/year in PPP])"
$ws.Range("B2").Value = 0.244781082578977
$ws.Range("C2").Value = 0.237490071485306
$ws.Range("D2").Value = 0.243781094527363
$ws.Range("E2").Value = 0.242201834862385
$ws.Range("F2").Value = 0.145118733509235
$ws.Range("G2").Value = 0.219512195121951
$ws.Range("H2").Value = 0.218354430379747
$ws.Range("I2").Value = 0.277777777777778
$ws.Range("J2").Value = 0.341880341880342
$ws.Range("K2").Value = 0.2
$ws.Range("L2").Value = 0.157446808510638
$ws.Range("M2").Value = 0.315463917525773

# Row 3
$ws.Range("A3").Value = "Supports tax on world top 3% to finance global poverty reduction`n(Additional 15% tax over [`$80k], 30% over [`$120k], 45% over [`$1M])"
$ws.Range("B3").Value = 0.288705924467514
$ws.Range("C3").Value = 0.295729250604351
$ws.Range("D3").Value = 0.252525252525253
$ws.Range("E3").Value = 0.312127236580517
$ws.Range("F3").Value = 0.254641909814324
$ws.Range("G3").Value = 0.240157480314961
$ws.Range("H3").Value = 0.289198606271777
$ws.Range("I3").Value = 0.267441860465116
$ws.Range("J3").Value = 0.519148936170213
$ws.Range("K3").Value = 0.284466019417476
$ws.Range("L3").Value = 0.145283018867925
$ws.Range("M3").Value = 0.329449838187702

# Row 4
$ws.Range("A4").Value = "Prefers sustainable future"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0

# Row 5
$ws.Range("A5").Value = "`"Governments should actively cooperate to have`nall countries converge in terms of GDP per capita by the end of the century`""
$ws.Range("B5").Value = 0.240272727272727
$ws.Range("C5").Value = 0.1972
$ws.Range("D5").Value = 0.18546365914787
$ws.Range("E5").Value = 0.21087786259542
$ws.Range("F5").Value = 0.113756613756614
$ws.Range("G5").Value = 0.142
$ws.Range("H5").Value = 0.139303482587065
$ws.Range("I5").Value = 0.288135593220339
$ws.Range("J5").Value = 0.294243070362473
$ws.Range("K5").Value = 0.242
$ws.Range("L5").Value = 0.063
$ws.Range("M5").Value = 0.37

# Row 6
$ws.Range("A6").Value = "Could sign a petition and spread ideas"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0

# Row 7
$ws.Range("A7").Value = "More likely to vote for party if part of worldwide`ncoalition for climate action and global redistribution"
$ws.Range("B7").Value = 0.1704
$ws.Range("C7").Value = 0.16
$ws.Range("D7").Value = 0.171679197994987
$ws.Range("E7").Value = 0.157442748091603
$ws.Range("F7").Value = 0.111111111111111
$ws.Range("G7").Value = 0.166
$ws.Range("H7").Value = 0.129353233830846
$ws.Range("I7").Value = 0.174334140435835
$ws.Range("J7").Value = 0.232409381663113
$ws.Range("K7").Value = 0.1765
$ws.Range("M7").Value = 0.183666666666667

# Row 8
$ws.Range("A8").Value = "Supports reparations for colonization and slavery in`nthe form of funding education and technology transfers"
$ws.Range("B8").Value = 0.414450291565922
$ws.Range("C8").Value = 0.384271892830563
$ws.Range("D8").Value = 0.422305764411028
$ws.Range("E8").Value = 0.425572519083969
$ws.Range("F8").Value = 0.247354497354497
$ws.Range("H8").Value = 0.386401326699834
$ws.Range("I8").Value = 0.418886198547215
$ws.Range("M8").Value = 0.455

# Row 9
$ws.Range("A9").Value = "`"My taxes should go towards solving global problems`""
$ws.Range("B9").Value = 0.266363636363636
$ws.Range("C9").Value = 0.2838
$ws.Range("D9").Value = 0.422305764411028
$ws.Range("E9").Value = 0.270038167938931
$ws.Range("F9").Value = 0.16005291005291
$ws.Range("G9").Value = 0.244
$ws.Range("H9").Value = 0.203980099502488
$ws.Range("I9").Value = 0.328087167070218
$ws.Range("J9").Value = 0.345415778251599
$ws.Range("K9").Value = 0.227
$ws.Range("L9").Value = 0.09
$ws.Range("M9").Value = 0.322333333333333
